$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Clear the obsolete employee phone/email cells ---
# Sarah (row 5): remove phone number
$ws.Range("F5").ClearContents()
# Nate (row 8): remove phone number
$ws.Range("F8").ClearContents()
# Janet (row 11): remove email and phone number
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()
# Bob (row 17): remove email (phone number F17 is kept)
$ws.Range("E17").ClearContents()

# --- Rebuild the mailto hyperlinks, dropping the ones for Janet (E11) and Bob (E17) ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:fred@email.com", "", "", "fred@email.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:hank@email.com", "", "", "hank@email.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:bill@email.com", "", "", "bill@email.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:sarah@email.com", "", "", "sarah@email.com")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:ted@email.com", "", "", "ted@email.com")
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:steve@email.com", "", "", "steve@email.com")
$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:nate@email.com", "", "", "nate@email.com")
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:jill@email.com", "", "", "jill@email.com")
$ws.Hyperlinks.Add($ws.Range("E10"), "mailto:carl@email.com", "", "", "carl@email.com")
$ws.Hyperlinks.Add($ws.Range("E12"), "mailto:clair@email.com", "", "", "clair@email.com")
$ws.Hyperlinks.Add($ws.Range("E13"), "mailto:justin@email.com", "", "", "justin@email.com")
$ws.Hyperlinks.Add($ws.Range("E14"), "mailto:tim@email.com", "", "", "tim@email.com")
$ws.Hyperlinks.Add($ws.Range("E15"), "mailto:victor@email.com", "", "", "victor@email.com")
$ws.Hyperlinks.Add($ws.Range("E16"), "mailto:samantha@email.com", "", "", "samantha@email.com")
$ws.Hyperlinks.Add($ws.Range("E18"), "mailto:adam@email.com", "", "", "adam@email.com")
$ws.Hyperlinks.Add($ws.Range("E19"), "mailto:denise@email.com", "", "", "denise@email.com")
$ws.Hyperlinks.Add($ws.Range("E20"), "mailto:gret@email.com", "", "", "gret@email.com")
$ws.Hyperlinks.Add($ws.Range("E21"), "mailto:lucy@email.com", "", "", "lucy@email.com")
$ws.Hyperlinks.Add($ws.Range("E22"), "mailto:megan@email.com", "", "", "megan@email.com")
$ws.Hyperlinks.Add($ws.Range("E23"), "mailto:sally@email.com", "", "", "sally@email.com")
$ws.Hyperlinks.Add($ws.Range("E24"), "mailto:frank@email.com", "", "", "frank@email.com")

# --- Leave the cursor where the author last edited ---
$ws.Range("F17").Select()
